$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new submission (Hossam Tabana / Medium Training) ---

# Numeric / date cells (A3 plain number, B3/C3 date-time matching row 2 style)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45192.45601851852
$ws.Range("C3").Value = 45192.45664351852
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat

# Text cells with real content
$ws.Range("D3").Value = 'Hossam.Ibrahim'
$ws.Range("E3").Value = 'Hossam Tabana'
$ws.Range("I3").Value = 'Geospatial Maps'
$ws.Range("L3").Value = 'Sometimes'
$ws.Range("O3").Value = 'Yes, simple calculations'
$ws.Range("R3").Value = 'Beginner level'
$ws.Range("U3").Value = 'Basic transformations only'
$ws.Range("X3").Value = 'Yes, basic automation'
$ws.Range("AA3").Value = 'Occasionally'
$ws.Range("AD3").Value = 'Very Important'
$ws.Range("AG3").Value = 'Yes, as static files'
$ws.Range("AJ3").Value = 'Occasionally'
$ws.Range("AM3").Value = 'Yes, it''s essential'
$ws.Range("AP3").Value = 'Not concerned'
$ws.Range("AS3").Value = 'Just exploring'
$ws.Range("AV3").Value = 'Yes, basic trend lines'
$ws.Range("AY3").Value = 'Beginner'
$ws.Range("BB3").Value = 'Definitely'
$ws.Range("BE3").Value = 'I know what it is but haven''t used it'
$ws.Range("BH3").Value = 'Possibly'
$ws.Range("BK3").Value = 'Yes, to multiple formats'
$ws.Range("BN3").Value = 'Using Power BI workspaces'
$ws.Range("BQ3").Value = 'Medium Training'

# Empty-string text cells (mirrors the blank "Points"/"Feedback" columns in row 2):
# a bare Value = "" clears the cell entirely, so force text-typed content via the
# apostrophe (quote-prefix) literal, then reset the style back to Normal so no
# stray quote-prefix formatting is left behind.
$emptyTextCols = @('F', 'G', 'H', 'J', 'K', 'M', 'N', 'P', 'Q', 'S', 'T', 'V', 'W', 'Y', 'Z', 'AB', 'AC', 'AE', 'AF', 'AH', 'AI', 'AK', 'AL', 'AN', 'AO', 'AQ', 'AR', 'AT', 'AU', 'AW', 'AX', 'AZ', 'BA', 'BC', 'BD', 'BF', 'BG', 'BI', 'BJ', 'BL', 'BM', 'BO', 'BP')
foreach ($col in $emptyTextCols) {
    $cell = $ws.Range("$col`3")
    $cell.Value = "'"
    $cell.Style = "Normal"
}
